$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.795.34"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "3.456.41"
$ws.Range("E3").Value = "  -3.82%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "597.48"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").Value = "137.75"
$ws.Range("E6").Value = "  -7.61%  "
$ws.Range("D7").Value = "3.453.62"
$ws.Range("E7").Value = "  -3.85%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").Value = "7.52"
$ws.Range("E10").Value = "  -5.63%  "
$ws.Range("E11").Value = "  -9.69%  "
$ws.Range("E12").Value = "  -7.85%  "
$ws.Range("D13").Value = "4.039.01"
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("E14").Value = "  -10.54%  "
$ws.Range("D15").Value = "26.71"
$ws.Range("E15").Value = "  -10.03%  "
$ws.Range("D16").Value = "3.456.98"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").Value = "65.709.66"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "0.115"
$ws.Range("E18").Value = "  -2.31%  "
$ws.Range("D19").Value = "9.93"
$ws.Range("E19").Value = "  -10.31%  "
$ws.Range("D20").Value = "5.81"
$ws.Range("E20").Value = "  -8.38%  "
$ws.Range("D21").Value = "'13.80"
$ws.Range("E21").Value = "  -7.25%  "
$ws.Range("D22").Value = "396.18"
$ws.Range("E22").Value = "  -6.74%  "
$ws.Range("D23").Value = "'0.550"
$ws.Range("E23").Value = "  -10.35%  "
$ws.Range("D24").Value = "73.63"
$ws.Range("E24").Value = "  -5.83%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "3.600.23"
$ws.Range("E26").Value = "  -3.58%  "
$ws.Range("D27").Value = "'0.0000108"
$ws.Range("E27").Value = "  -9.93%  "
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").Value = "7.31"
$ws.Range("E29").Value = "  -10.25%  "
$ws.Range("E30").Value = "  -8.69%  "
$ws.Range("D31").Value = "8.27"
$ws.Range("E31").Value = "  -12.02%  "
$ws.Range("D32").Value = "3.459.15"
$ws.Range("E32").Value = "  -3.57%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -7.74%  "
$ws.Range("D35").Value = "'23.00"
$ws.Range("E35").Value = "  -8.18%  "
$ws.Range("D36").Value = "173.58"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("E37").Value = "  -14.15%  "
$ws.Range("D38").Value = "6.98"
$ws.Range("E38").Value = "  -9.95%  "
$ws.Range("D39").Value = "1.54"
$ws.Range("E39").Value = "  -7.68%  "
$ws.Range("D40").Value = "4.85"
$ws.Range("E40").Value = "  -12.65%  "
$ws.Range("D41").Value = "0.0783"
$ws.Range("E41").Value = "  -8.33%  "
$ws.Range("D42").Value = "0.822"
$ws.Range("E42").Value = "  -6.75%  "
$ws.Range("E43").Value = "  -5.51%  "
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "4.46"
$ws.Range("E45").Value = "  -14.16%  "
$ws.Range("E46").Value = "  -11.44%  "
$ws.Range("D47").Value = "23.49"
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("D48").Value = "1.12"
$ws.Range("E48").Value = "  -1.77%  "
$ws.Range("D49").Value = "'6.60"
$ws.Range("E49").Value = "  -7.57%  "
$ws.Range("D50").Value = "2.14"
$ws.Range("E50").Value = "  -15.35%  "
$ws.Range("D51").Value = "2.209.39"
$ws.Range("E51").Value = "  -8.47%  "
